$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16 contains a table (graphicFrame "Google Shape;213;p29") whose
#    table style is switched from the deck's custom style to the built-in
#    "Medium Style 2 - Accent 1" gallery style, exactly like picking a new
#    style from the Table Design > Table Styles gallery in the ribbon.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item("Google Shape;213;p29")
$table = $tableShape.Table
$table.ApplyStyle("{8869D2DC-A362-4009-9C6E-DE220DEC3937}")

# ---------------------------------------------------------------------------
# 2) The presentation's theme colour palette is switched from the custom
#    "Integral" palette to the stock "Office" palette (Design > Variants >
#    Colors > "Office"), i.e. every theme colour slot gets the standard
#    Office RGB values.
# ---------------------------------------------------------------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
$colorScheme.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # Dark 1    000000
$colorScheme.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # Light 1   FFFFFF
$colorScheme.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # Dark 2    44546A
$colorScheme.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # Light 2   E7E6E6
$colorScheme.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # Accent 1  5B9BD5
$colorScheme.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # Accent 2  ED7D31
$colorScheme.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # Accent 3  A5A5A5
$colorScheme.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # Accent 4  FFC000
$colorScheme.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # Accent 5  4472C4
$colorScheme.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # Accent 6  70AD47
$colorScheme.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # Hyperlink 0563C1
$colorScheme.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # Followed Hyperlink 954F72
